$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.207.14"
$ws.Range("E2").Value = "  +1.75%  "

$ws.Range("D3").Value = "3.349.60"
$ws.Range("E3").Value = "  +2.21%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "191.20"
$ws.Range("E5").Value = "  +2.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "591.27"
$ws.Range("E6").Value = "  +1.57%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.605"
$ws.Range("E8").Value = "  +0.53%  "

$ws.Range("E9").Value = "  +1.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.76"
$ws.Range("E10").Value = "  +3.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.417"
$ws.Range("E11").Value = "  +1.25%  "

$ws.Range("D12").Value = "3.943.10"
$ws.Range("E12").Value = "  +2.57%  "

$ws.Range("E13").Value = "  -1.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.43"
$ws.Range("E14").Value = "  +2.93%  "

$ws.Range("D15").Value = "69.328.41"
$ws.Range("E15").Value = "  +1.92%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000170"
$ws.Range("E16").Value = "  +0.88%  "

$ws.Range("D17").Value = "3.355.62"
$ws.Range("E17").Value = "  +2.30%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "448.09"
$ws.Range("E18").Value = "  +12.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.81"
$ws.Range("E19").Value = "  +1.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.75"
$ws.Range("E20").Value = "  +1.18%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.81"
$ws.Range("E21").Value = "  +2.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.94"
$ws.Range("E22").Value = "  +6.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.21%  "

$ws.Range("B24").Value = "WrappedeETH"
$ws.Range("C24").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D24").Value = "3.514.04"
$ws.Range("E24").Value = "  +2.84%  "

$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.523"
$ws.Range("E25").Value = "  +1.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000121"
$ws.Range("E26").Value = "  +2.63%  "

$ws.Range("E27").Value = "  +1.54%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.42"
$ws.Range("E28").Value = "  -1.25%  "

$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.01"
$ws.Range("E30").Value = "  +2.57%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.31"
$ws.Range("E31").Value = "  +2.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.54"
$ws.Range("E32").Value = "  +0.26%  "

$ws.Range("E33").Value = "  +1.92%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.98"
$ws.Range("E34").Value = "  +0.11%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.56"
$ws.Range("E36").Value = "  +5.99%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.96"
$ws.Range("E37").Value = "  +0.69%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.93"
$ws.Range("E38").Value = "  +1.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.21"
$ws.Range("E39").Value = "  +1.25%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.60"
$ws.Range("E40").Value = "  +1.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.805"
$ws.Range("E41").Value = "  -0.72%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.47"
$ws.Range("E42").Value = "  -0.96%  "

$ws.Range("D43").Value = "2.704.42"
$ws.Range("E43").Value = "  +0.91%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.50"
$ws.Range("E44").Value = "  +2.15%  "

$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.20"
$ws.Range("E45").Value = "  +1.11%  "

$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0689"
$ws.Range("E46").Value = "  +0.35%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.36"
$ws.Range("E47").Value = "  +2.16%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "334.54"
$ws.Range("E48").Value = "  -0.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0284"
$ws.Range("E49").Value = "  +2.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "32.44"
$ws.Range("E50").Value = "  +5.17%  "

$ws.Range("E51").Value = "  +3.30%  "
